# Rotate the "Meeting Year" ... "Census Division 9" block of rows in
# Table 1 up by one position (the label/coefficient/std-error of each
# row moves to the row above it; the former top row's content wraps
# around to the bottom), per the commit "update substantive regressions
# & preliminaries".
#
# Table layout (1-based Word row numbers):
#   13/14  Meeting Year                    | 0.006***  / (0.001)
#   15/16  Overall Climate Vulnerability   | -0.132    / (0.115)
#   17/18  Census Division 2               | -0.015    / (0.017)
#   19/20  Census Division 3               | -0.027    / (0.017)
#   21/22  Census Division 4               | -0.036**  / (0.018)
#   23/24  Census Division 5               | -0.011    / (0.018)
#   25/26  Census Division 6               | -0.020    / (0.018)
#   27/28  Census Division 7               | -0.006    / (0.019)
#   29/30  Census Division 8               | -0.011    / (0.019)
#   31/32  Census Division 9               | 0.039*    / (0.023)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each entry: label row, value/std row, new label, new coefficient, new std error
$newRows = @(
    @(13, 14, "Overall Climate Vulnerability", "-0.132",   "(0.115)"),
    @(15, 16, "Census Division 2",              "-0.015",   "(0.017)"),
    @(17, 18, "Census Division 3",              "-0.027",   "(0.017)"),
    @(19, 20, "Census Division 4",              "-0.036**", "(0.018)"),
    @(21, 22, "Census Division 5",              "-0.011",   "(0.018)"),
    @(23, 24, "Census Division 6",              "-0.020",   "(0.018)"),
    @(25, 26, "Census Division 7",              "-0.006",   "(0.019)"),
    @(27, 28, "Census Division 8",              "-0.011",   "(0.019)"),
    @(29, 30, "Census Division 9",              "0.039*",   "(0.023)"),
    @(31, 32, "Meeting Year",                   "0.006***", "(0.001)")
)

foreach ($entry in $newRows) {
    $labelRow = $entry[0]
    $stdRow   = $entry[1]
    $label    = $entry[2]
    $coef     = $entry[3]
    $stderr   = $entry[4]

    $t.Cell($labelRow, 1).Range.Text = $label
    $t.Cell($labelRow, 2).Range.Text = $coef
    $t.Cell($stdRow, 2).Range.Text   = $stderr
}
